$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1) and "全部类型" (sheet4) both contain a duplicated
# table of convention info; update the "想去人数" (F column) figures in
# each to reflect the newly scraped counts.

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("F2").Value = 699
    $ws.Range("F3").Value = 23
    $ws.Range("F4").Value = 533
    $ws.Range("F8").Value = 48
    $ws.Range("F9").Value = 4435
    $ws.Range("F10").Value = 4322
    $ws.Range("F11").Value = 13
    $ws.Range("F12").Value = 133
}
